$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.040552626876211
$ws.Range("D2").Value = 1.043197025379953
$ws.Range("E2").Value = 1.047581029945516
$ws.Range("F2").Value = 1.055499912855133
$ws.Range("I2").Value = 1.027521332220243
$ws.Range("J2").Value = 1.045638614573486
$ws.Range("K2").Value = 1.045971570969731
$ws.Range("L2").Value = 1.050343268142906
$ws.Range("M2").Value = 1.058240200855259
$ws.Range("N2").Value = 1.018718237031132

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.042677659721558
$ws.Range("D3").Value = 1.045236962835933
$ws.Range("E3").Value = 1.04945847728362
$ws.Range("F3").Value = 1.057433623175898
$ws.Range("I3").Value = 1.027548965694502
$ws.Range("J3").Value = 1.047403452191046
$ws.Range("K3").Value = 1.047819283169714
$ws.Range("L3").Value = 1.052029828701801
$ws.Range("M3").Value = 1.059984510230903
$ws.Range("N3").Value = 1.019338789505845

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.044047690269625
$ws.Range("D4").Value = 1.046552335988866
$ws.Range("E4").Value = 1.050668597753903
$ws.Range("F4").Value = 1.058679433074639
$ws.Range("I4").Value = 1.027564016835793
$ws.Range("J4").Value = 1.048540331842387
$ws.Range("K4").Value = 1.049009904582399
$ws.Range("L4").Value = 1.053116024412633
$ws.Range("M4").Value = 1.061107365112031
$ws.Range("N4").Value = 1.019737683351249

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.044622481656577
$ws.Range("D5").Value = 1.047104244066611
$ws.Range("E5").Value = 1.051176228544913
$ws.Range("F5").Value = 1.059201896416888
$ws.Range("I5").Value = 1.027569667625648
$ws.Range("J5").Value = 1.049017082341346
$ws.Range("K5").Value = 1.049509276821298
$ws.Range("L5").Value = 1.053571457918981
$ws.Range("M5").Value = 1.061578040742703
$ws.Range("N5").Value = 1.019904753669177

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.044718923932512
$ws.Range("D6").Value = 1.047196849618628
$ws.Range("E6").Value = 1.05126139796773
$ws.Range("F6").Value = 1.059289546253037
$ws.Range("I6").Value = 1.02757057675512
$ws.Range("J6").Value = 1.04909706155469
$ws.Range("K6").Value = 1.049593056031967
$ws.Range("L6").Value = 1.053647857318362
$ws.Range("M6").Value = 1.061656989432815
$ws.Range("N6").Value = 1.019932769128691

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.044055375215838
$ws.Range("D7").Value = 1.046559714798297
$ws.Range("E7").Value = 1.050675385045649
$ws.Range("F7").Value = 1.058686419232321
$ws.Range("I7").Value = 1.027564094999896
$ws.Range("J7").Value = 1.048546706863407
$ws.Range("K7").Value = 1.049016581762295
$ws.Range("L7").Value = 1.053122114634079
$ws.Range("M7").Value = 1.061113659664831
$ws.Range("N7").Value = 1.019739918195873

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041271846395986
$ws.Range("D8").Value = 1.043887401260763
$ws.Range("E8").Value = 1.048216512907943
$ws.Range("F8").Value = 1.05615455865839
$ws.Range("I8").Value = 1.027531257391537
$ws.Range("J8").Value = 1.046236118083355
$ws.Range("K8").Value = 1.04659705875029
$ws.Range("L8").Value = 1.050914323599226
$ws.Range("M8").Value = 1.058830919519351
$ws.Range("N8").Value = 1.018928508536267

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036327227002844
$ws.Range("D9").Value = 1.039141963374086
$ws.Range("E9").Value = 1.043846453447684
$ws.Range("F9").Value = 1.051650385318968
$ws.Range("I9").Value = 1.027451689026105
$ws.Range("J9").Value = 1.042124487954971
$ws.Range("K9").Value = 1.04229432561427
$ws.Range("L9").Value = 1.046983645758231
$ws.Range("M9").Value = 1.054762740730566
$ws.Range("N9").Value = 1.017478053235575

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.033002253735928
$ws.Range("D10").Value = 1.035952079611793
$ws.Range("E10").Value = 1.040906491269352
$ws.Range("F10").Value = 1.048617299438485
$ws.Range("I10").Value = 1.027384000459417
$ws.Range("J10").Value = 1.039354893896583
$ws.Range("K10").Value = 1.039397849408819
$ws.Range("L10").Value = 1.04433465063845
$ws.Range("M10").Value = 1.052018415414991
$ws.Range("N10").Value = 1.016496650368783

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.031555307326735
$ws.Range("D11").Value = 1.034564207205569
$ws.Range("E11").Value = 1.039626795410642
$ws.Range("F11").Value = 1.047296394770215
$ws.Range("I11").Value = 1.02735120511151
$ws.Range("J11").Value = 1.038148516231383
$ws.Range("K11").Value = 1.038136641657674
$ws.Range("L11").Value = 1.043180505700038
$ws.Range("M11").Value = 1.050822119756158
$ws.Range("N11").Value = 1.016068144778474

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.031016726993757
$ws.Range("D12").Value = 1.034047659464746
$ws.Range("E12").Value = 1.039150426190296
$ws.Range("F12").Value = 1.046804585308307
$ws.Range("I12").Value = 1.027338498581138
$ws.Range("J12").Value = 1.037699312858788
$ws.Range("K12").Value = 1.03766708754779
$ws.Range("L12").Value = 1.042750707549676
$ws.Range("M12").Value = 1.050376533637235
$ws.Range("N12").Value = 1.015908434312794

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.031132305580853
$ws.Range("D13").Value = 1.034158507872891
$ws.Range("E13").Value = 1.039252656257609
$ws.Range("F13").Value = 1.046910133405504
$ws.Range("I13").Value = 1.027341247946384
$ws.Range("J13").Value = 1.037795718824701
$ws.Range("K13").Value = 1.037767858126908
$ws.Range("L13").Value = 1.04284295083514
$ws.Range("M13").Value = 1.05047216944256
$ws.Range("N13").Value = 1.015942717561092

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.03151081117043
$ws.Range("D14").Value = 1.034521530420164
$ws.Range("E14").Value = 1.039587439838616
$ws.Range("F14").Value = 1.047255765627807
$ws.Range("I14").Value = 1.027350165500123
$ws.Range("J14").Value = 1.038111407578827
$ws.Range("K14").Value = 1.038097850494686
$ws.Range("L14").Value = 1.043145001004168
$ws.Range("M14").Value = 1.050785312702376
$ws.Range("N14").Value = 1.016054954230589

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.031743871437572
$ws.Range("D15").Value = 1.034745063069175
$ws.Range("E15").Value = 1.039793573216302
$ws.Range("F15").Value = 1.047468565515191
$ws.Range("I15").Value = 1.027355590311617
$ws.Range("J15").Value = 1.03830576702505
$ws.Range("K15").Value = 1.03830102490495
$ws.Range("L15").Value = 1.043330957806629
$ws.Range("M15").Value = 1.050978086914896
$ws.Range("N15").Value = 1.016124034480617

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.033098127762292
$ws.Range("D16").Value = 1.036044045538463
$ws.Range("E16").Value = 1.040991277272795
$ws.Range("F16").Value = 1.048704801641955
$ws.Range("I16").Value = 1.027386103438012
$ws.Range("J16").Value = 1.03943480449685
$ws.Range("K16").Value = 1.039481401145958
$ws.Range("L16").Value = 1.044411095140432
$ws.Range("M16").Value = 1.052097638922656
$ws.Range("N16").Value = 1.016525013125345

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.033945660703903
$ws.Range("D17").Value = 1.036857063748003
$ws.Range("E17").Value = 1.041740757378222
$ws.Range("F17").Value = 1.049478214511313
$ws.Range("I17").Value = 1.027404309461618
$ws.Range("J17").Value = 1.040141090881504
$ws.Range("K17").Value = 1.040219919900547
$ws.Range("L17").Value = 1.045086712372716
$ws.Range("M17").Value = 1.05279774624877
$ws.Range("N17").Value = 1.016775577978096

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034439319072968
$ws.Range("D18").Value = 1.037330645963545
$ws.Range("E18").Value = 1.042177274709807
$ws.Range("F18").Value = 1.049928605996445
$ws.Range("I18").Value = 1.027414592485711
$ws.Range("J18").Value = 1.040552370392338
$ws.Range("K18").Value = 1.040650010761193
$ws.Range("L18").Value = 1.045480104179902
$ws.Range("M18").Value = 1.053205338587194
$ws.Range("N18").Value = 1.016921386272358

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034607527330816
$ws.Range("D19").Value = 1.037492018135359
$ws.Range("E19").Value = 1.042326007910901
$ws.Range("F19").Value = 1.050082055447934
$ws.Range("I19").Value = 1.027418041733671
$ws.Range("J19").Value = 1.040692490739381
$ws.Range("K19").Value = 1.04079654714756
$ws.Range("L19").Value = 1.045614125414607
$ws.Range("M19").Value = 1.053344187575249
$ws.Range("N19").Value = 1.016971045457183

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033854800410446
$ws.Range("D20").Value = 1.036769900776375
$ws.Range("E20").Value = 1.041660411854093
$ws.Range("F20").Value = 1.049395310066655
$ws.Range("I20").Value = 1.027402390915085
$ws.Range("J20").Value = 1.040065384149003
$ws.Range("K20").Value = 1.040140753834941
$ws.Range("L20").Value = 1.045014296024126
$ws.Range("M20").Value = 1.052722711048523
$ws.Range("N20").Value = 1.016748730187451

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.031399381923311
$ws.Range("D21").Value = 1.034414658039047
$ws.Range("E21").Value = 1.039488883188596
$ws.Range("F21").Value = 1.0471540180106
$ws.Range("I21").Value = 1.027347554001759
$ws.Range("J21").Value = 1.038018475749514
$ws.Range("K21").Value = 1.038000706185486
$ws.Range("L21").Value = 1.043056085282801
$ws.Range("M21").Value = 1.050693133969558
$ws.Range("N21").Value = 1.016021918444989

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.029849064233168
$ws.Range("D22").Value = 1.032927846259191
$ws.Range("E22").Value = 1.038117563577863
$ws.Range("F22").Value = 1.04573806305028
$ws.Range("I22").Value = 1.02731003845339
$ws.Range("J22").Value = 1.036725117694935
$ws.Range("K22").Value = 1.036648877507409
$ws.Range("L22").Value = 1.041818517048812
$ws.Range("M22").Value = 1.049409934989711
$ws.Range("N22").Value = 1.015561787773161

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.030671545359475
$ws.Range("D23").Value = 1.033716611315359
$ws.Range("E23").Value = 1.038845104528753
$ws.Range("F23").Value = 1.046489339390022
$ws.Range("I23").Value = 1.027330214474175
$ws.Range("J23").Value = 1.037411366866077
$ws.Range("K23").Value = 1.0373661149083
$ws.Range("L23").Value = 1.04247518828931
$ws.Range("M23").Value = 1.050090868081296
$ws.Range("N23").Value = 1.015806014491281

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.033895858406796
$ws.Range("D24").Value = 1.036809287940161
$ws.Range("E24").Value = 1.041696718519882
$ws.Range("F24").Value = 1.049432773262468
$ws.Range("I24").Value = 1.027403258862695
$ws.Range("J24").Value = 1.040099594879382
$ws.Range("K24").Value = 1.04017652765578
$ws.Range("L24").Value = 1.04504701996836
$ws.Range("M24").Value = 1.05275661859908
$ws.Range("N24").Value = 1.016760862604633

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037610424135921
$ws.Range("D25").Value = 1.040373272461677
$ws.Range("E25").Value = 1.044980788503476
$ws.Range("F25").Value = 1.052820045408856
$ws.Range("I25").Value = 1.027474836432801
$ws.Range("J25").Value = 1.043192351513035
$ws.Range("K25").Value = 1.043411498271781
$ws.Range("L25").Value = 1.048004741973763
$ws.Range("M25").Value = 1.055820028967252
$ws.Range("N25").Value = 1.017855533650782
